$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("P2").ClearContents()
$ws.Range("Q2").ClearContents()
$ws.Range("R2").ClearContents()
$ws.Range("S2").ClearContents()
$ws.Range("T2").ClearContents()
$ws.Range("U2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("W2").ClearContents()
$ws.Range("X2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("AC2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()
$ws.Range("AF2").ClearContents()
$ws.Range("AI2").ClearContents()
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AJ2").Value = 78294468

# Row 3
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("P3").ClearContents()
$ws.Range("Q3").ClearContents()
$ws.Range("R3").ClearContents()
$ws.Range("S3").ClearContents()
$ws.Range("T3").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("V3").ClearContents()
$ws.Range("W3").ClearContents()
$ws.Range("X3").ClearContents()
$ws.Range("Y3").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AA3").ClearContents()
$ws.Range("AB3").ClearContents()
$ws.Range("AC3").ClearContents()
$ws.Range("AD3").ClearContents()
$ws.Range("AE3").ClearContents()
$ws.Range("AF3").ClearContents()
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 82331835

# Row 4
$ws.Range("Y4").ClearContents()
$ws.Range("Z4").ClearContents()
$ws.Range("D4").Value = 2172
$ws.Range("E4").Value = -97
$ws.Range("F4").Value = -97
$ws.Range("G4").Value = 48
$ws.Range("H4").Value = 32
$ws.Range("I4").Value = 33
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 6940
$ws.Range("L4").Value = 5717
$ws.Range("M4").Value = 1223
$ws.Range("N4").Value = 1223
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 866
$ws.Range("Q4").Value = 620
$ws.Range("R4").Value = -266
$ws.Range("S4").Value = -320
$ws.Range("T4").Value = 206
$ws.Range("U4").Value = 414
$ws.Range("V4").Value = 1981
$ws.Range("W4").Value = -4.49
$ws.Range("X4").Value = 1.46
$ws.Range("AA4").Value = 467.56
$ws.Range("AB4").Value = 118.95
$ws.Range("AC4").Value = 34
$ws.Range("AD4").Value = 65.16
$ws.Range("AE4").Value = 707
$ws.Range("AF4").Value = 3.13
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 172993713

# Row 5
$ws.Range("D5").Value = 5745
$ws.Range("E5").Value = -58
$ws.Range("F5").Value = -58
$ws.Range("G5").Value = -149
$ws.Range("H5").Value = -392
$ws.Range("I5").Value = -391
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 6829
$ws.Range("L5").Value = 6010
$ws.Range("M5").Value = 818
$ws.Range("N5").Value = 819
$ws.Range("O5").Value = -1
$ws.Range("P5").Value = 866
$ws.Range("Q5").Value = -171
$ws.Range("R5").Value = -220
$ws.Range("S5").Value = 132
$ws.Range("T5").Value = 183
$ws.Range("U5").Value = -354
$ws.Range("V5").Value = 2409
$ws.Range("W5").Value = -1
$ws.Range("X5").Value = -6.82
$ws.Range("Y5").Value = -38.33
$ws.Range("Z5").Value = -5.69
$ws.Range("AA5").Value = 734.6799999999999
$ws.Range("AB5").Value = 76.48
$ws.Range("AC5").Value = -226
$ws.Range("AD5").Value = -8.470000000000001
$ws.Range("AE5").Value = 697
$ws.Range("AF5").Value = 2.75
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 172993713

# Row 6
$ws.Range("D6").Value = 4247
$ws.Range("E6").Value = -42
$ws.Range("F6").Value = -42
$ws.Range("G6").Value = -204
$ws.Range("H6").Value = 107
$ws.Range("I6").Value = 107
$ws.Range("K6").Value = 4131
$ws.Range("L6").Value = 2991
$ws.Range("M6").Value = 1141
$ws.Range("N6").Value = 1141
$ws.Range("P6").Value = 866
$ws.Range("Q6").Value = -445
$ws.Range("R6").Value = 808
$ws.Range("S6").Value = -587
$ws.Range("T6").Value = 114
$ws.Range("U6").Value = -559
$ws.Range("V6").Value = 1771
$ws.Range("W6").Value = -0.98
$ws.Range("X6").Value = 2.52
$ws.Range("Y6").Value = 10.92
$ws.Range("Z6").Value = 1.95
$ws.Range("AA6").Value = 262.14
$ws.Range("AB6").Value = 107.77
$ws.Range("AC6").Value = 62
$ws.Range("AD6").Value = 16.25
$ws.Range("AE6").Value = 805
$ws.Range("AF6").Value = 1.25
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 172993713

# Row 7
$ws.Range("U7").ClearContents()
$ws.Range("AI7").ClearContents()
$ws.Range("D7").Value = 4570
$ws.Range("E7").Value = 132
$ws.Range("G7").Value = 30
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 4
$ws.Range("K7").Value = 5180
$ws.Range("L7").Value = 3990
$ws.Range("M7").Value = 1190
$ws.Range("N7").Value = 1190
$ws.Range("P7").Value = 870
$ws.Range("Q7").Value = 40
$ws.Range("R7").Value = -180
$ws.Range("S7").Value = 150
$ws.Range("T7").Value = 130
$ws.Range("W7").Value = 2.89
$ws.Range("X7").Value = 0.09
$ws.Range("Y7").Value = 0.34
$ws.Range("Z7").Value = 0.09
$ws.Range("AA7").Value = 335.29
$ws.Range("AC7").Value = 2
$ws.Range("AD7").Value = 385.67
$ws.Range("AE7").Value = 835
$ws.Range("AF7").Value = 1.06
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0

# Row 8
$ws.Range("U8").ClearContents()
$ws.Range("AI8").ClearContents()
$ws.Range("D8").Value = 5163
$ws.Range("E8").Value = 281
$ws.Range("G8").Value = 210
$ws.Range("H8").Value = 166
$ws.Range("I8").Value = 166
$ws.Range("K8").Value = 5320
$ws.Range("L8").Value = 3960
$ws.Range("M8").Value = 1360
$ws.Range("N8").Value = 1360
$ws.Range("P8").Value = 870
$ws.Range("Q8").Value = 290
$ws.Range("R8").Value = -210
$ws.Range("S8").Value = -10
$ws.Range("T8").Value = 110
$ws.Range("W8").Value = 5.44
$ws.Range("X8").Value = 3.21
$ws.Range("Y8").Value = 13.02
$ws.Range("Z8").Value = 3.16
$ws.Range("AA8").Value = 291.18
$ws.Range("AC8").Value = 96
$ws.Range("AD8").Value = 9.300000000000001
$ws.Range("AE8").Value = 954
$ws.Range("AF8").Value = 0.93
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0

# Row 9
$ws.Range("U9").ClearContents()
$ws.Range("AI9").ClearContents()
$ws.Range("D9").Value = 5390
$ws.Range("E9").Value = 300
$ws.Range("G9").Value = 240
$ws.Range("H9").Value = 190
$ws.Range("I9").Value = 190
$ws.Range("K9").Value = 5460
$ws.Range("L9").Value = 3920
$ws.Range("M9").Value = 1540
$ws.Range("N9").Value = 1540
$ws.Range("P9").Value = 870
$ws.Range("Q9").Value = 310
$ws.Range("R9").Value = -220
$ws.Range("S9").Value = -40
$ws.Range("T9").Value = 110
$ws.Range("W9").Value = 5.57
$ws.Range("X9").Value = 3.52
$ws.Range("Y9").Value = 13.1
$ws.Range("Z9").Value = 3.52
$ws.Range("AA9").Value = 254.55
$ws.Range("AC9").Value = 109
$ws.Range("AD9").Value = 8.130000000000001
$ws.Range("AE9").Value = 1081
$ws.Range("AF9").Value = 0.82
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
